$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 308
$ws1.Range("F3").Value = 68
$ws1.Range("F4").Value = 3722
$ws1.Range("F5").Value = 2261
$ws1.Range("F8").Value = 11
$ws1.Range("F9").Value = 181
$ws1.Range("F10").Value = 101
$ws1.Range("F11").Value = 82
$ws1.Range("F12").Value = 1378
$ws1.Range("F14").Value = 2181

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 34
$ws2.Range("F3").Value = 2

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 308
$ws4.Range("F3").Value = 68
$ws4.Range("F4").Value = 3722
$ws4.Range("F5").Value = 2261
$ws4.Range("F8").Value = 11
$ws4.Range("F9").Value = 34
$ws4.Range("F10").Value = 181
$ws4.Range("F11").Value = 101
$ws4.Range("F12").Value = 82
$ws4.Range("F13").Value = 2
$ws4.Range("F15").Value = 1378
$ws4.Range("F17").Value = 2181
